# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), as three new, separate runs:
#   <w:r><w:t xml:space="preserve"> (</w:t></w:r>
#   <w:r><w:t>Changed main</w:t></w:r>
#   <w:r><w:t>)</w:t></w:r>
# while leaving the original run ("This is a Microsoft word document.")
# completely untouched.
#
# A plain Range.InsertAfter() would just append the new text into the
# existing run (because the new text would share the same - empty -
# run formatting), merging everything into a single <w:r>. To force
# Word to create genuinely separate run elements we replace the
# matched range using Range.InsertXML() with a flat WordprocessingML
# "xmlPackage" document fragment that already contains the desired run
# boundaries; InsertXML never coalesces the runs it inserts with their
# neighbours.

$d = $word.ActiveDocument

$original = "This is a Microsoft word document."

# Locate the sentence we need to extend.
$target = $d.Content
$found = $target.Find.Execute($original, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to edit."
}

$originalXml = $original.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
$suffixXml = '<w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>Changed main</w:t></w:r><w:r><w:t>)</w:t></w:r>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r><w:t>' + $originalXml + '</w:t></w:r>' +
    $suffixXml +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Replace the matched range (the whole sentence) with itself plus the
# three new trailing runs. Because InsertXML substitutes the range's
# contents wholesale, the original run is re-created verbatim and the
# new runs stay distinct siblings instead of being merged into it.
$target.InsertXML($packageXml)
